$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.203.97'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.60%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.822.08'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.38'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.15%  '

$ws.Range("E6").Value = '  -1.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4250'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.52%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3680'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07237'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.87%  '

$ws.Range("E10").Value = '  -3.17%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.98'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.87%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.820.20'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.26%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.701'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07093'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.307'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.70%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.10'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008871'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.006'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.96%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.05'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.232.40'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.49%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.137'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.84%  '

$ws.Range("E23").Value = '  -3.54%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.047.78'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.80%  '

$ws.Range("E25").Value = '  -1.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.38'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.35'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.122'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.99%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.215'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -4.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.11'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.73%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08865'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.93%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.194'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7546'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.428'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.56%  '

$ws.Range("E35").Value = '  -6.43%  '

$ws.Range("E36").Value = '  -1.19%  '

$ws.Range("E37").Value = '  -2.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01971'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.50%  '

$ws.Range("E39").Value = '  -1.19%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.135'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.73%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.869'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.72%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1690'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.32%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5035'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.53%  '

$ws.Range("E44").Value = '  -1.97%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.62'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.49%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '106.95'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.46%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4736'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.55%  '

$ws.Range("E48").Value = '  -1.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06371'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.659'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.53%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.807'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.32%  '
